# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 05:31"

# Kazajistan (row 31) - updated case counts
$ws.Range("B31").Value = 103815
$ws.Range("C31").Value = 244
$ws.Range("D31").Value = 86450
$ws.Range("E31").Value = 15950

# Belgica (row 40) - updated case counts
$ws.Range("B40").Value = 79479
$ws.Range("C40").Value = 582
$ws.Range("D40").Value = 18078
$ws.Range("E40").Value = 51432
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 9969

# Honduras (row 51) - updated case counts
$ws.Range("B51").Value = 52298
$ws.Range("C51").Value = 628
$ws.Range("D51").Value = 7867
$ws.Range("E51").Value = 42823
$ws.Range("G51").Value = 15
$ws.Range("H51").Value = 1608

# Australia / Austria swap order (rows 71-72) with updated Australia counts
$ws.Range("A71").Value = "Australia"
$ws.Range("B71").Value = 24236
$ws.Range("C71").Value = 243
$ws.Range("D71").Value = 15248
$ws.Range("E71").Value = 8525
$ws.Range("G71").Value = 13
$ws.Range("H71").Value = 463

$ws.Range("A72").Value = "Austria"
$ws.Range("B72").Value = 24084
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 20958
$ws.Range("E72").Value = 2397
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 729

# Camboya (row 183) - updated case counts
$ws.Range("D183").Value = 253
$ws.Range("E183").Value = 20

# Islas Malvinas / Montserrat swap order (rows 213-214)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
